$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 79 (which currently holds FSZZ130),
# shifting FSZZ130 down to row 81.
$ws.Rows.Item(79).Resize(2).Insert()

# Fill in the new rows with the added records.
$ws.Range("A79").Value = "FSZZ106"
$ws.Range("B79").Value = 0.3112873072124995
$ws.Range("C79").Value = 0.6327834472285486

$ws.Range("A80").Value = "FSZZ108"
$ws.Range("B80").Value = 0.170393098269208
$ws.Range("C80").Value = 0.5191002556620498
